# 25A element - changed dV from 5% to 10%
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("relays")

# Update "Max. Slip Voltage [%]" (column C) from 5 to 10 for the affected rows
$rows = @(2,3,4,5,6,7,8,9,11,12,13)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = 10
}

# Update the active cell selection to match the edited sheet state
$ws.Range("D11").Select()
